$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the CasesTab query (row 2, column B) to remove the Cohort clause
$newQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
"MATCH (c)<--(diag:diagnosis)`n" +
"MATCH (samp:sample)-->(c) `n" +
"  MATCH (f:file)-[*]->(c)`n" +
"   WHERE f.file_type IN [`"RNA Sequence File`"] `n" +
"OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
"  WITH DISTINCT c, s, demo, diag, co`n" +
"RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
"        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
"        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
"        coalesce(demo.breed, '') AS Breed ,`n" +
"        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
"        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
"        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
"        coalesce(demo.sex, '') AS Sex ,`n" +
"        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
"        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
"        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newQuery

# Update sheet view: clear the scrolled top-left cell and move selection to B2
$ws.Activate()
$ws.Range("B2").Select()

$wb.Save()
